$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed rows 14:15 with the formatting of the last existing data row (row 13)
# so the new rows reuse the same style indices instead of minting new ones.
$ws.Range("A13:F13").Copy()
$ws.Range("A14:F15").PasteSpecial()

# New row 14: 2025-10-07 (serial 45937), 四方坪站
$ws.Range("A14").Value = 45937
$ws.Range("B14").Value = "四方坪站"
$ws.Range("C14").Value = 10341.38
$ws.Range("D14").Value = 8749.08
$ws.Range("E14").Value = 3562.56
$ws.Range("F14").Value = 396

# New row 15: 2025-10-07 (serial 45937), 高岭站
$ws.Range("A15").Value = 45937
$ws.Range("B15").Value = "高岭站"
$ws.Range("C15").Value = 4056.06
$ws.Range("D15").Value = 3260.87
$ws.Range("E15").Value = 1094.62
$ws.Range("F15").Value = 143

# Update the active selection to H11 as recorded in the saved view state
$ws.Range("H11").Select()
